$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Rename the header row: "_old" -> "_FV2210" and "_new" -> "_FV2304"
#    suffixes (columns A:J hold the "_old"/FV2210 block, K holds the
#    constant "diff" header, columns L:U hold the "_new"/FV2304 block).
# ---------------------------------------------------------------------------
$fv2210Headers = @(
    "Segmentname_FV2210",
    "Segmentgruppe_FV2210",
    "Segment_FV2210",
    "Datenelement_FV2210",
    "Segment ID_FV2210",
    "Code_FV2210",
    "Qualifier_FV2210",
    "Beschreibung_FV2210",
    "Bedingungsausdruck_FV2210",
    "Bedingung_FV2210"
)

$fv2304Headers = @(
    "Segmentname_FV2304",
    "Segmentgruppe_FV2304",
    "Segment_FV2304",
    "Datenelement_FV2304",
    "Segment ID_FV2304",
    "Code_FV2304",
    "Qualifier_FV2304",
    "Beschreibung_FV2304",
    "Bedingungsausdruck_FV2304",
    "Bedingung_FV2304"
)

for ($i = 0; $i -lt 10; $i++) {
    $ws.Cells.Item(1, $i + 1).Value2 = $fv2210Headers[$i]
}
for ($i = 0; $i -lt 10; $i++) {
    $ws.Cells.Item(1, $i + 12).Value2 = $fv2304Headers[$i]
}

# ---------------------------------------------------------------------------
# 2. Turn the used range into an Excel Table ("Table1") so the new header
#    names drive the table column headers as well.
# ---------------------------------------------------------------------------
$tableRange = $ws.Range("A1:U92")
$lo = $ws.ListObjects.Add(1, $tableRange, 0, 1)
$lo.Name = "Table1"
$lo.TableStyle = ""

# ---------------------------------------------------------------------------
# 3. Freeze the header row (split below row 1, frozen, active pane bottom
#    left, starting at A2).
# ---------------------------------------------------------------------------
$ws.Range("A2").Select() | Out-Null
$excel.ActiveWindow.FreezePanes = $true
